# Update the city list: several cities in each state were swapped out for a
# new set of cities (used to drive the new "selected state" image + graph
# connector UI). Column A (state) is untouched; only specific column B
# (city) cells change.
#
# The write order below matters: it reproduces the exact shared-string
# table layout of the saved workbook (new strings are appended to the
# shared-string table in the order the underlying cells are written).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colorado block (rows 2-11) - replace 5 cities, written in descending row order
$ws.Range("B11").Value = "Aspen"
$ws.Range("B10").Value = "Burlington"
$ws.Range("B9").Value = "Gunnison"
$ws.Range("B7").Value = "Durango"
$ws.Range("B6").Value = "Springfield"

# New York block (rows 22-31) - replace 5 cities, written in ascending row order
$ws.Range("B23").Value = "Dunkirk"
$ws.Range("B25").Value = "Watertown"
$ws.Range("B27").Value = "Binghamton"
$ws.Range("B28").Value = "Albany"
$ws.Range("B29").Value = "Kingston"

# California block (rows 12-21) - replace 3 cities, written in ascending row order
$ws.Range("B18").Value = "Redding"
$ws.Range("B20").Value = "Eureka"
$ws.Range("B21").Value = "Santa Barbara"

# Update the saved selection/view to cell C21 (was B30, with topLeftCell A21)
[void]$ws.Range("C21").Select()
